# Edit script: update the "K" column (column G) values on Sheet1
# per regenerated save_data (commit: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals").
#
# Only the K (strikeouts) values in column G change for this sheet; all
# other columns / cells are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new K value (column G)
$kValues = @{
    2 = 1
    3 = 0
    4 = 2
    5 = 1
    6 = 1
    7 = 2
    8 = 1
    9 = 2
    10 = 0
    11 = 2
    12 = 1
    13 = 2
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 2
    25 = 1
    26 = 2
    27 = 1
    28 = 1
    29 = 3
    30 = 1
    31 = 1
    32 = 2
    34 = 1
    35 = 0
    36 = 2
    37 = 0
    38 = 1
    39 = 0
    40 = 0
    41 = 1
    42 = 0
    43 = 1
    44 = 2
    45 = 1
    46 = 0
    47 = 0
    48 = 1
    49 = 1
    50 = 2
    51 = 1
    52 = 0
    53 = 0
    54 = 0
    55 = 2
    56 = 1
    57 = 1
    58 = 1
    59 = 0
    60 = 1
    61 = 1
    62 = 1
    63 = 1
    64 = 2
    65 = 2
    66 = 1
    67 = 1
    68 = 2
    69 = 2
    70 = 1
    71 = 2
    72 = 1
    75 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}
